$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Fermentation TAL yield (row 23): Midpoint E23 0.68 -> 0.73
$ws.Range("E23").Value = 0.73

# Update Fermentation TAL titer (row 24): Midpoint E24 76 -> 68
$ws.Range("E24").Value = 68

# Update TAL decarboxylation conversion (row 28): Midpoint E28 0.05 -> 0.048
# Also change G28/I28 formulas to reference E28 instead of hardcoded 0.05
$ws.Range("E28").Value = 0.048
$ws.Range("G28").Formula = "=E28*0.0463/0.2087"
$ws.Range("I28").Formula = "=E28*0.34/0.2087"

# Update the selection on the sheet to A28:XFD29 (activeCell A28)
$ws.Range("A28:XFD29").Select()

$wb.Save()
